# Update the Testing Report with current Pass/Fail results and contributors
# for the requirements that have now been tested (current to 27/03).
#
# Column D = "Pass or Fail" (P: Pass, F: Fail)
# Column E = "Contributor" (R: Rommel, J: Jaidyn)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Requirements 1-5 (rows 6-10): log in, log out, make account,
# delete account, reset password -> Passed, tested by Rommel
$passRows = 6, 7, 8, 9, 10

# Requirement 11 + its sub-items (rows 16-20): change profile
# (changing PFP/name/short bio, roles) -> Failed, tested by Jaidyn
$failRows = 16, 17, 18, 19, 20

# Fill in column D (Pass or Fail) first, for all rows
foreach ($r in $passRows) {
    $ws.Cells.Item($r, 4).Value = "P"
}
foreach ($r in $failRows) {
    $ws.Cells.Item($r, 4).Value = "F"
}

# Then fill in column E (Contributor)
foreach ($r in $passRows) {
    $ws.Cells.Item($r, 5).Value = "R"
}
foreach ($r in $failRows) {
    $ws.Cells.Item($r, 5).Value = "J"
}

# Leave the cursor/selection where the author left off working (row 11 area)
$ws.Range("C11").Select()
